# Correcting bug with spaces in filenames
#
# 1) One sheet's chart has a bogus category-axis title -- a Python list
#    literal that leaked into the title text instead of being used as
#    actual category labels: ["NUEDEXTA User", "Non User"]. Clear it; the
#    axis keeps its <c:title> element (layout/overlay) but loses the
#    <c:tx> rich text run.
# 2) Every sheet's C3 header cell changes from "T3" to " T2" (leading
#    space) -- the actual filename/label bug fix referenced by the
#    commit message.

$wb = $excel.ActiveWorkbook

$badTitle = '["NUEDEXTA User", "Non User"]'

foreach ($ws in $wb.Worksheets) {

    # --- 1) Strip the bogus category-axis title wherever it shows up ---
    $chartCount = $ws.ChartObjects().Count
    for ($i = 1; $i -le $chartCount; $i++) {
        $chart = $ws.ChartObjects().Item($i).Chart
        $catAx = $chart.Axes(1)   # xlCategory
        if ($catAx.HasTitle -and $catAx.AxisTitle.Text -eq $badTitle) {
            $catAx.AxisTitle.Text = ""
            $catAx.HasTitle = $false
        }
    }

    # --- 2) Fix the " T2" header text ---
    $ws.Range("C3").Value = " T2"
}
